$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 - header row. Keep existing Tool/Disk labels, rename Time column
# into two columns (Download time / Build time), shift Disk/Version and
# add a new Notes column.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Download time"
$ws.Range("C1").Value = "Build time"
$ws.Range("D1").Value = "Disk"
$ws.Range("E1").Value = "Version"
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Value = "Notes"
$ws.Range("F1").Font.Bold = $true

# ---------------------------------------------------------------------
# Row 2 - DRAM
# ---------------------------------------------------------------------
$oldC2 = $ws.Range("C2").Value2
$oldD2 = $ws.Range("D2").Value2
$ws.Range("D2").Value = $oldC2
$ws.Range("E2").Value = $oldD2
$ws.Range("C2").Value = ""
$ws.Range("F2").Value = "This is what the IT department told us"

# ---------------------------------------------------------------------
# Row 3 - InterProScan
# ---------------------------------------------------------------------
$oldC3 = $ws.Range("C3").Value2
$oldD3 = $ws.Range("D3").Value2
$ws.Range("D3").Value = $oldC3
$ws.Range("E3").Value = $oldD3
$ws.Range("B3").Value = "~ 2:34:27"
$ws.Range("B3").NumberFormat = "hh:mm:ss\ AM/PM"
$ws.Range("C3").Value = "~ 11:32"

# ---------------------------------------------------------------------
# Row 4 - Metacerberus
# ---------------------------------------------------------------------
$oldC4 = $ws.Range("C4").Value2
$oldD4 = $ws.Range("D4").Value2
$ws.Range("D4").Value = $oldC4
$ws.Range("E4").Value = $oldD4
$ws.Range("C4").Value = 0

# ---------------------------------------------------------------------
# Row 5 - PROKKA
# ---------------------------------------------------------------------
$oldC5 = $ws.Range("C5").Value2
$oldD5 = $ws.Range("D5").Value2
$ws.Range("D5").Value = $oldC5
$ws.Range("E5").Value = $oldD5
$ws.Range("C5").Value = 0
$ws.Range("F5").Value = "Everything is downloaded together with bioconda install"

# ---------------------------------------------------------------------
# Row 6 - eggNOG (new row)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "eggNOG"
$ws.Range("D6").Value = "31GB"
$ws.Range("F6").Value = "Includes archaea, bacteria, and virus HMMs"

# ---------------------------------------------------------------------
# Row 7 - Microbeannotator (new row)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "Microbeannotator"

# ---------------------------------------------------------------------
# Row 33 - move the trailing "." marker from G33 to H33
# ---------------------------------------------------------------------
$ws.Range("G33").Clear()
$ws.Range("H33").Value = "."

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.33
$ws.Columns.Item(2).ColumnWidth = 13.65
$ws.Columns.Item(6).ColumnWidth = 43.65

# ---------------------------------------------------------------------
# Selection, matching the saved cursor position of the edited workbook
# ---------------------------------------------------------------------
[void]$ws.Range("G6").Select()
